$wb = $excel.ActiveWorkbook

# --- Sheet: Overall ---
# A2 was a numeric 3018; convert to text "3,018"
$ws1 = $wb.Worksheets.Item("Overall")
$ws1.Range("A2").NumberFormat = "@"
$ws1.Range("A2").Value = "3,018"

# --- Sheet: County ---
$ws2 = $wb.Worksheets.Item("County")

# Existing B column counts converted from numeric to text (same displayed value)
$ws2.Range("B2:B15").NumberFormat = "@"
$ws2.Range("B2").Value = "151"
$ws2.Range("B3").Value = "119"
$ws2.Range("B4").Value = "147"
$ws2.Range("B5").Value = "30"
$ws2.Range("B6").Value = "271"
$ws2.Range("B7").Value = "35"
$ws2.Range("B8").Value = "153"
$ws2.Range("B9").Value = "96"
$ws2.Range("B10").Value = "702"
$ws2.Range("B11").Value = "31"
$ws2.Range("B12").Value = "222"
$ws2.Range("B13").Value = "118"
$ws2.Range("B14").Value = "718"
$ws2.Range("B15").Value = "225"

# New "Total" row appended as row 16
$ws2.Range("A16:F16").NumberFormat = "@"
$ws2.Range("A16").Value = "Total"
$ws2.Range("B16").Value = "3,018"
$ws2.Range("C16").Value = "$10,165,747,853"
$ws2.Range("D16").Value = "11.74%"
$ws2.Range("E16").Value = "-5.38%"
$ws2.Range("F16").Value = "60.30%"

# --- Sheet: Congressional District ---
$ws3 = $wb.Worksheets.Item("Congressional District")
$ws3.Range("B2:B11").NumberFormat = "@"
$ws3.Range("B2").Value = "309"
$ws3.Range("B3").Value = "286"
$ws3.Range("B4").Value = "230"
$ws3.Range("B5").Value = "267"
$ws3.Range("B6").Value = "357"
$ws3.Range("B7").Value = "237"
$ws3.Range("B8").Value = "539"
$ws3.Range("B9").Value = "434"
$ws3.Range("B10").Value = "359"
$ws3.Range("B11").Value = "3,018"

# --- Sheet: Size ---
$ws4 = $wb.Worksheets.Item("Size")
$ws4.Range("B2:B8").NumberFormat = "@"
$ws4.Range("B2").Value = "772"
$ws4.Range("B3").Value = "928"
$ws4.Range("B4").Value = "489"
$ws4.Range("B5").Value = "230"
$ws4.Range("B6").Value = "483"
$ws4.Range("B7").Value = "116"
$ws4.Range("B8").Value = "3,018"

# --- Sheet: Subsector ---
$ws5 = $wb.Worksheets.Item("Subsector")
$ws5.Range("B2:B14").NumberFormat = "@"
$ws5.Range("B2").Value = "361"
$ws5.Range("B3").Value = "372"
$ws5.Range("B4").Value = "125"
$ws5.Range("B5").Value = "234"
$ws5.Range("B6").Value = "8"
$ws5.Range("B7").Value = "852"
$ws5.Range("B8").Value = "45"
$ws5.Range("B9").Value = "1"
$ws5.Range("B10").Value = "242"
$ws5.Range("B11").Value = "35"
$ws5.Range("B12").Value = "687"
$ws5.Range("B13").Value = "56"
$ws5.Range("B14").Value = "3,018"
